$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Col1a1"
$ws.Cells.Item(2, 3).Value = "Ddr2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 9.108069666666667
$ws.Cells.Item(2, 8).Value = 27.324209
$ws.Cells.Item(2, 9).Value = 0.00155006418458712
$ws.Cells.Item(2, 10).Value = 0.00155006418458712
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 1.4763
$ws.Cells.Item(2, 14).Value = 4.428900000000001
$ws.Cells.Item(2, 15).Value = 0.01318769285519422
$ws.Cells.Item(2, 16).Value = 0.01318769285519422
$ws.Cells.Item(2, 17).Value = 13.4462432489
$ws.Cells.Item(2, 18).Value = 121.0161892401
$ws.Cells.Item(2, 19).Value = 0.00002044177037217201
$ws.Cells.Item(2, 20).Value = 0.00002044177037217201

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Col1a1"
$ws.Cells.Item(3, 3).Value = "Ddr2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 9.108069666666667
$ws.Cells.Item(3, 8).Value = 27.324209
$ws.Cells.Item(3, 9).Value = 0.00155006418458712
$ws.Cells.Item(3, 10).Value = 0.00155006418458712
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 78.17189533333334
$ws.Cells.Item(3, 14).Value = 234.515686
$ws.Cells.Item(3, 15).Value = 0.6983045082736506
$ws.Cells.Item(3, 16).Value = 0.6983045082736506
$ws.Cells.Item(3, 17).Value = 711.995068671375
$ws.Cells.Item(3, 18).Value = 6407.955618042374
$ws.Cells.Item(3, 19).Value = 0.001082416808210706
$ws.Cells.Item(3, 20).Value = 0.001082416808210706

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Col1a1"
$ws.Cells.Item(4, 3).Value = "Ddr2"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 9.108069666666667
$ws.Cells.Item(4, 8).Value = 27.324209
$ws.Cells.Item(4, 9).Value = 0.00155006418458712
$ws.Cells.Item(4, 10).Value = 0.00155006418458712
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.1536526666666667
$ws.Cells.Item(4, 14).Value = 0.460958
$ws.Cells.Item(4, 15).Value = 0.00137256937911098
$ws.Cells.Item(4, 16).Value = 0.00137256937911098
$ws.Cells.Item(4, 17).Value = 1.399479192469111
$ws.Cells.Item(4, 18).Value = 12.595312732222
$ws.Cells.Item(4, 19).Value = 0.000002127570635420909
$ws.Cells.Item(4, 20).Value = 0.000002127570635420909

$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Col1a1"
$ws.Cells.Item(5, 3).Value = "Ddr2"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 9.108069666666667
$ws.Cells.Item(5, 8).Value = 27.324209
$ws.Cells.Item(5, 9).Value = 0.00155006418458712
$ws.Cells.Item(5, 10).Value = 0.00155006418458712
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 32.14343433333334
$ws.Cells.Item(5, 14).Value = 96.43030300000001
$ws.Cells.Item(5, 15).Value = 0.2871352294920441
$ws.Cells.Item(5, 16).Value = 0.2871352294920441
$ws.Cells.Item(5, 17).Value = 292.7646392339253
$ws.Cells.Item(5, 18).Value = 2634.881753105327
$ws.Cells.Item(5, 19).Value = 0.0004450780353688207
$ws.Cells.Item(5, 20).Value = 0.0004450780353688207

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Col1a1"
$ws.Cells.Item(6, 3).Value = "Ddr2"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 5771.873535333333
$ws.Cells.Item(6, 8).Value = 17315.620606
$ws.Cells.Item(6, 9).Value = 0.9822909543423312
$ws.Cells.Item(6, 10).Value = 0.9822909543423313
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 1.4763
$ws.Cells.Item(6, 14).Value = 4.428900000000001
$ws.Cells.Item(6, 15).Value = 0.01318769285519422
$ws.Cells.Item(6, 16).Value = 0.01318769285519422
$ws.Cells.Item(6, 17).Value = 8521.0169002126
$ws.Cells.Item(6, 18).Value = 76689.15210191341
$ws.Cells.Item(6, 19).Value = 0.01295415140030227
$ws.Cells.Item(6, 20).Value = 0.01295415140030228

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Col1a1"
$ws.Cells.Item(7, 3).Value = "Ddr2"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 5771.873535333333
$ws.Cells.Item(7, 8).Value = 17315.620606
$ws.Cells.Item(7, 9).Value = 0.9822909543423312
$ws.Cells.Item(7, 10).Value = 0.9822909543423313
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 78.17189533333334
$ws.Cells.Item(7, 14).Value = 234.515686
$ws.Cells.Item(7, 15).Value = 0.6983045082736506
$ws.Cells.Item(7, 16).Value = 0.6983045082736506
$ws.Cells.Item(7, 17).Value = 451198.293881314
$ws.Cells.Item(7, 18).Value = 4060784.644931826
$ws.Cells.Item(7, 19).Value = 0.6859382018536766
$ws.Cells.Item(7, 20).Value = 0.6859382018536767

$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Col1a1"
$ws.Cells.Item(8, 3).Value = "Ddr2"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 5771.873535333333
$ws.Cells.Item(8, 8).Value = 17315.620606
$ws.Cells.Item(8, 9).Value = 0.9822909543423312
$ws.Cells.Item(8, 10).Value = 0.9822909543423313
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.1536526666666667
$ws.Cells.Item(8, 14).Value = 0.460958
$ws.Cells.Item(8, 15).Value = 0.00137256937911098
$ws.Cells.Item(8, 16).Value = 0.00137256937911098
$ws.Cells.Item(8, 17).Value = 886.8637603667274
$ws.Cells.Item(8, 18).Value = 7981.773843300548
$ws.Cells.Item(8, 19).Value = 0.001348262485307985
$ws.Cells.Item(8, 20).Value = 0.001348262485307985

$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Col1a1"
$ws.Cells.Item(9, 3).Value = "Ddr2"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 5771.873535333333
$ws.Cells.Item(9, 8).Value = 17315.620606
$ws.Cells.Item(9, 9).Value = 0.9822909543423312
$ws.Cells.Item(9, 10).Value = 0.9822909543423313
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 32.14343433333334
$ws.Cells.Item(9, 14).Value = 96.43030300000001
$ws.Cells.Item(9, 15).Value = 0.2871352294920441
$ws.Cells.Item(9, 16).Value = 0.2871352294920441
$ws.Cells.Item(9, 17).Value = 185527.8379632915
$ws.Cells.Item(9, 18).Value = 1669750.541669624
$ws.Cells.Item(9, 19).Value = 0.2820503386030442
$ws.Cells.Item(9, 20).Value = 0.2820503386030443

$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Col1a1"
$ws.Cells.Item(10, 3).Value = "Ddr2"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.272029666666667
$ws.Cells.Item(10, 8).Value = 3.816089
$ws.Cells.Item(10, 9).Value = 0.0002164813950916887
$ws.Cells.Item(10, 10).Value = 0.0002164813950916887
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 1.4763
$ws.Cells.Item(10, 14).Value = 4.428900000000001
$ws.Cells.Item(10, 15).Value = 0.01318769285519422
$ws.Cells.Item(10, 16).Value = 0.01318769285519422
$ws.Cells.Item(10, 17).Value = 1.8778973969
$ws.Cells.Item(10, 18).Value = 16.9010765721
$ws.Cells.Item(10, 19).Value = 0.000002854890147333141
$ws.Cells.Item(10, 20).Value = 0.000002854890147333141

$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Col1a1"
$ws.Cells.Item(11, 3).Value = "Ddr2"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 1.272029666666667
$ws.Cells.Item(11, 8).Value = 3.816089
$ws.Cells.Item(11, 9).Value = 0.0002164813950916887
$ws.Cells.Item(11, 10).Value = 0.0002164813950916887
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 78.17189533333334
$ws.Cells.Item(11, 14).Value = 234.515686
$ws.Cells.Item(11, 15).Value = 0.6983045082736506
$ws.Cells.Item(11, 16).Value = 0.6983045082736506
$ws.Cells.Item(11, 17).Value = 99.43696996356155
$ws.Cells.Item(11, 18).Value = 894.932729672054
$ws.Cells.Item(11, 19).Value = 0.0001511699341498956
$ws.Cells.Item(11, 20).Value = 0.0001511699341498956

$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Col1a1"
$ws.Cells.Item(12, 3).Value = "Ddr2"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 1.272029666666667
$ws.Cells.Item(12, 8).Value = 3.816089
$ws.Cells.Item(12, 9).Value = 0.0002164813950916887
$ws.Cells.Item(12, 10).Value = 0.0002164813950916887
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 0.6666666666666666
$ws.Cells.Item(12, 13).Value = 0.1536526666666667
$ws.Cells.Item(12, 14).Value = 0.460958
$ws.Cells.Item(12, 15).Value = 0.00137256937911098
$ws.Cells.Item(12, 16).Value = 0.00137256937911098
$ws.Cells.Item(12, 17).Value = 0.1954507503624444
$ws.Cells.Item(12, 18).Value = 1.759056753262
$ws.Cells.Item(12, 19).Value = 0.0000002971357340500778
$ws.Cells.Item(12, 20).Value = 0.0000002971357340500778

$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Col1a1"
$ws.Cells.Item(13, 3).Value = "Ddr2"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 1.272029666666667
$ws.Cells.Item(13, 8).Value = 3.816089
$ws.Cells.Item(13, 9).Value = 0.0002164813950916887
$ws.Cells.Item(13, 10).Value = 0.0002164813950916887
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 32.14343433333334
$ws.Cells.Item(13, 14).Value = 96.43030300000001
$ws.Cells.Item(13, 15).Value = 0.2871352294920441
$ws.Cells.Item(13, 16).Value = 0.2871352294920441
$ws.Cells.Item(13, 17).Value = 40.88740206055189
$ws.Cells.Item(13, 18).Value = 367.986618544967
$ws.Cells.Item(13, 19).Value = 0.0000621594350604099
$ws.Cells.Item(13, 20).Value = 0.0000621594350604099

$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Col1a1"
$ws.Cells.Item(14, 3).Value = "Ddr2"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 93.67702500000001
$ws.Cells.Item(14, 8).Value = 281.031075
$ws.Cells.Item(14, 9).Value = 0.01594250007799006
$ws.Cells.Item(14, 10).Value = 0.01594250007799006
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 0.6666666666666666
$ws.Cells.Item(14, 13).Value = 1.4763
$ws.Cells.Item(14, 14).Value = 4.428900000000001
$ws.Cells.Item(14, 15).Value = 0.01318769285519422
$ws.Cells.Item(14, 16).Value = 0.01318769285519422
$ws.Cells.Item(14, 17).Value = 138.2953920075
$ws.Cells.Item(14, 18).Value = 1244.6585280675
$ws.Cells.Item(14, 19).Value = 0.0002102447943724429
$ws.Cells.Item(14, 20).Value = 0.0002102447943724429

$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Col1a1"
$ws.Cells.Item(15, 3).Value = "Ddr2"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 93.67702500000001
$ws.Cells.Item(15, 8).Value = 281.031075
$ws.Cells.Item(15, 9).Value = 0.01594250007799006
$ws.Cells.Item(15, 10).Value = 0.01594250007799006
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 78.17189533333334
$ws.Cells.Item(15, 14).Value = 234.515686
$ws.Cells.Item(15, 15).Value = 0.6983045082736506
$ws.Cells.Item(15, 16).Value = 0.6983045082736506
$ws.Cells.Item(15, 17).Value = 7322.910593438051
$ws.Cells.Item(15, 18).Value = 65906.19534094246
$ws.Cells.Item(15, 19).Value = 0.01113271967761349
$ws.Cells.Item(15, 20).Value = 0.01113271967761349

$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Col1a1"
$ws.Cells.Item(16, 3).Value = "Ddr2"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 93.67702500000001
$ws.Cells.Item(16, 8).Value = 281.031075
$ws.Cells.Item(16, 9).Value = 0.01594250007799006
$ws.Cells.Item(16, 10).Value = 0.01594250007799006
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 0.1536526666666667
$ws.Cells.Item(16, 14).Value = 0.460958
$ws.Cells.Item(16, 15).Value = 0.00137256937911098
$ws.Cells.Item(16, 16).Value = 0.00137256937911098
$ws.Cells.Item(16, 17).Value = 14.39372469665
$ws.Cells.Item(16, 18).Value = 129.54352226985
$ws.Cells.Item(16, 19).Value = 0.00002188218743352356
$ws.Cells.Item(16, 20).Value = 0.00002188218743352356

$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Col1a1"
$ws.Cells.Item(17, 3).Value = "Ddr2"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 93.67702500000001
$ws.Cells.Item(17, 8).Value = 281.031075
$ws.Cells.Item(17, 9).Value = 0.01594250007799006
$ws.Cells.Item(17, 10).Value = 0.01594250007799006
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 32.14343433333334
$ws.Cells.Item(17, 14).Value = 96.43030300000001
$ws.Cells.Item(17, 15).Value = 0.2871352294920441
$ws.Cells.Item(17, 16).Value = 0.2871352294920441
$ws.Cells.Item(17, 17).Value = 3011.101301629526
$ws.Cells.Item(17, 18).Value = 27099.91171466573
$ws.Cells.Item(17, 19).Value = 0.004577653418570606
$ws.Cells.Item(17, 20).Value = 0.004577653418570606

